# Updated symbol list on Fri Feb  3 20:51:59 UTC 2023 with GitHub Actions
# Refresh crypto price/volume snapshot values (columns D=Price, E=Volume(1h))
# and fix the ordering of the BOLO / CoinbaseStockToken rows (48 & 49).
# NumberFormat is forced to "@" (Text) before writing each D/E cell so that
# Excel does not auto-convert these numeric-looking strings (e.g. "330.50",
# "0.40%") into real numbers/percentages, which would lose their original
# textual formatting (trailing zeros, "%" suffix, etc.).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '330.50'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.40%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '41.29'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '1.95%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.683'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.95%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08063'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.67%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.019'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '2.54%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.750'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.28%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '4.524'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-1.77%'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.73%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9241'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-2.20%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1260'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-3.25%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1945'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-2.34%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.284'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-7.15%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09307'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.23%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.03704'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '6.09%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.1053'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '9.46%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.001300'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.17%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006221'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.03%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.25%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-2.53%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1418'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.29%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2655'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '10.02%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04421'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.38%'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.03%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004348'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.22%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001243'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '13.89%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02846'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '15.94%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05464'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2.92%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007755'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '3.58%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009970'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '15.07%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1421'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-0.86%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002135'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '4.10%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01172'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '11.49%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006770'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-2.04%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.00%'
$ws.Range("B48").Value = 'BOLO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003013'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-14.16%'
$ws.Range("B49").Value = 'CoinbaseStockToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002283'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '34.06%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002104'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.00%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002004'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.00%'
